# Round 24 tip updates - update match day dates in column A
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDates = @{
    2  = 45885
    3  = 45884
    4  = 45885
    5  = 45885
    6  = 45884
    7  = 45884
    8  = 45885
    9  = 45886
    10 = 45885
    11 = 45885
    12 = 45885
    13 = 45886
    14 = 45885
    15 = 45886
    16 = 45884
    17 = 45886
    18 = 45886
    19 = 45886
}

foreach ($row in $newDates.Keys) {
    $ws.Cells.Item($row, 1).Value = $newDates[$row]
}

$wb.Save()
